$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("وسپهر", "سرمایه گذاریها", 56069000000, 1.5561182349, 0.1819371276, 0.01769802, 2.84611454)
    ,@("وبصادر", "بانک ها، موسسات اعتباری و سایر نهادهای پولی", 48833099640, 1.5014039866, 0.1453584896, 0.01593854, 2.47881351)
    ,@("تیپیکو", "مواد و محصولات دارویی", 70100000000, 0.7326509394, 0.0784467814, 0.01141072, 3.55834114)
    ,@("وبملت", "بانک ها، موسسات اعتباری و سایر نهادهای پولی", 58068000000, 0.8856630173, 0.1058060463, 0.03747018, 2.94758564)
    ,@("فملی", "فلزات اساسی", 48300000000, 0.823199619, 0.0900707746, 0.0237642, 2.45175288)
    ,@("تاصیکو", "استخراج کانه های فلزی", 31775000000, 0.3539362408, 0.0582551103, 0.02728387, 1.61292853)
    ,@("سفارس", "سیمان، آهک و گچ", 72762500000, 0.9621651542, 0.1077471765, 0.04516604, 3.69349211)
    ,@("انتخاب", "ماشین آلات و تجهیزات", 1016494856340, 1.2927220328, 0.1345797348, 0.01170711, 51.59822346)
    ,@("شستا", "شرکتهای چند رشته ای صنعتی", 52147200000, 1.1114061145, 0.1040262314, 0.02129395, 2.64704033)
    ,@("وکار", "بانک ها، موسسات اعتباری و سایر نهادهای پولی", 246819441838, 0.4363923146, 0.0775113587, 0.00328698, 12.52878422)
    ,@("فباهنر", "فلزات اساسی", 15750000000, 1.2967432778, 0.1546501591, 0.03468104, 0.79948464)
    ,@("فولاد", "فلزات اساسی", 252900000000, 0.8213573532, 0.0843446747, 0.01697103, 12.837439)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}